# Apply cryptos list price/volume refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "27.844.75"
Set-TextValue "E2" "  +2.43%  "
Set-TextValue "D3" "1.872.52"
Set-TextValue "E3" "  +0.78%  "
Set-TextValue "D4" "1.013"
Set-TextValue "D5" "313.35"
Set-TextValue "E5" "  +0.65%  "
Set-TextValue "E6" "  -0.62%  "
Set-TextValue "D7" "0.4830"
Set-TextValue "E7" "  +0.69%  "
Set-TextValue "D8" "0.3818"
Set-TextValue "E8" "  +2.77%  "
Set-TextValue "E9" "  +1.02%  "
Set-TextValue "D10" "0.9398"
Set-TextValue "E10" "  +0.52%  "
Set-TextValue "E11" "  +4.67%  "
Set-TextValue "E12" "  -0.95%  "
Set-TextValue "D13" "1.889.73"
Set-TextValue "E13" "  +1.55%  "
Set-TextValue "D14" "5.510"
Set-TextValue "D15" "6.610"
Set-TextValue "E15" "  +1.14%  "
Set-TextValue "D16" "91.31"
Set-TextValue "E16" "  +1.50%  "
Set-TextValue "E17" "  -0.62%  "
Set-TextValue "D18" "0.000008845"
Set-TextValue "E18" "  +1.35%  "
Set-TextValue "E19" "  -0.67%  "
Set-TextValue "D20" "27.877.80"
Set-TextValue "E20" "  +2.39%  "
Set-TextValue "D21" "14.85"
Set-TextValue "E21" "  +1.11%  "
Set-TextValue "E22" "  +0.56%  "
Set-TextValue "D23" "2.118.45"
Set-TextValue "E23" "  +0.90%  "
Set-TextValue "D24" "10.86"
Set-TextValue "D25" "157.73"
Set-TextValue "E25" "  +2.60%  "
Set-TextValue "D26" "1.947"
Set-TextValue "E26" "  -0.37%  "
Set-TextValue "E27" "  +0.54%  "
Set-TextValue "D28" "2.044"
Set-TextValue "E28" "  +2.52%  "
Set-TextValue "D29" "116.01"
Set-TextValue "E29" "  +0.29%  "
Set-TextValue "D30" "4.976"
Set-TextValue "E30" "  +0.85%  "
Set-TextValue "D31" "0.08891"
Set-TextValue "E31" "  -0.01%  "
Set-TextValue "D32" "3.339"
Set-TextValue "E32" "  +0.73%  "
Set-TextValue "D33" "1.223"
Set-TextValue "E33" "  +3.41%  "
Set-TextValue "E34" "  +4.74%  "
Set-TextValue "D35" "4.654"
Set-TextValue "E35" "  +1.58%  "
Set-TextValue "D36" "2.726"
Set-TextValue "E36" "  +1.67%  "
Set-TextValue "D37" "1.130"
Set-TextValue "E37" "  +0.67%  "
Set-TextValue "D38" "0.02043"
Set-TextValue "E38" "  +1.49%  "
Set-TextValue "E39" "  +5.38%  "
Set-TextValue "D40" "0.05375"
Set-TextValue "E40" "  +2.45%  "
Set-TextValue "D41" "3.002"
Set-TextValue "E41" "  +0.23%  "
Set-TextValue "D42" "7.058"
Set-TextValue "E42" "  -0.26%  "
Set-TextValue "D43" "8.534"
Set-TextValue "E43" "  +2.38%  "
Set-TextValue "D44" "0.1530"
Set-TextValue "E44" "  +0.23%  "
Set-TextValue "E45" "  +0.66%  "
Set-TextValue "D46" "0.4871"
Set-TextValue "E46" "  +2.12%  "
Set-TextValue "D47" "105.57"
Set-TextValue "E47" "  +3.05%  "
Set-TextValue "E48" "  -0.66%  "
Set-TextValue "D49" "1.664"
Set-TextValue "E49" "  +2.07%  "
Set-TextValue "D50" "68.19"
Set-TextValue "E50" "  +2.75%  "
Set-TextValue "D51" "0.06124"
Set-TextValue "E51" "  +0.79%  "
